$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Connectors")

# --- Row 11: G-Switch DC-058A-D010 barrel jack ---------------------------
$ws.Range("A11").Value = "DC-058A-D010"
$ws.Range("B11").Value = "JACK `nDC-058A-D010 TH  HORIZ"
$ws.Range("C11").Value = "DC-058A-D010"
$ws.Range("D11").Value = "TH"
$ws.Range("E11").Value = "NA"
$ws.Range("F11").Value = "12V"
$ws.Range("G11").Value = "5A"
$ws.Range("H11").Value = "HORZ 2WIRE Jack"
$ws.Range("I11").Value = "NA"
$ws.Range("J11").Value = "G-Switch"
$ws.Range("K11").Value = "CON_DC-058A-D010"
$ws.Range("L11").Value = "Altium_Footprints.PcbLib"
$ws.Range("M11").Value = "JACK_2P"
$ws.Range("N11").Value = "Altium_Schematic_Symbols.SchLib"
$ws.Range("O11").Value = "Datasheet"
$ws.Range("P11").Value = "https://datasheet.lcsc.com/lcsc/2102241737_G-Switch-DC-058A-D010_C2686970.pdf"

# --- Row 12: Molex 436500200 receptacle -----------------------------------
$ws.Range("A12").Value = 436500200
$ws.Range("B12").Value = "RECEP 436500200 TH 3mm HORIZ"
$ws.Range("C12").Value = 436500200
$ws.Range("D12").Value = "TH"
$ws.Range("E12").Value = "NA"
$ws.Range("G12").Value = "8.5A"
$ws.Range("F12").Value = "600V"
$ws.Range("H12").Value = "RERCEP 2WIRE HORIZ"
$ws.Range("I12").Value = "NA"
$ws.Range("J12").Value = "MOLEX"
$ws.Range("K12").Value = "CON_436500200"
$ws.Range("L12").Value = "Altium_Footprints.PcbLib"
$ws.Range("M12").Value = "JACK_2P"
$ws.Range("N12").Value = "Altium_Schematic_Symbols.SchLib"
$ws.Range("O12").Value = "Datasheet"
$ws.Range("P12").Value = "https://datasheet.lcsc.com/lcsc/2305301754_MOLEX-436500200_C192562.pdf"

# --- Formatting: both new rows mirror the wrap-text style used elsewhere --
$ws.Range("A11:P11").WrapText = $true
$ws.Range("A12:P12").WrapText = $true

# Match the auto-computed row heights Excel produced for these rows.
$ws.Rows.Item(11).RowHeight = 60
$ws.Rows.Item(12).RowHeight = 30

# --- View state: Connectors tab becomes the active / selected sheet -------
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 3
$ws.Range("P14").Select()
